$wb = $excel.ActiveWorkbook

# Rename the "Include from " sheet to "Include from Diagnostic Evide"
$wsInclude = $wb.Worksheets.Item("Include from ")
$wsInclude.Name = "Include from Diagnostic Evide"

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "0.0.0"
$wsMeta.Range("B5").Value = "Diagnostic Evidence Base"
$wsMeta.Range("B7").Formula = "=""false"""
$wsMeta.Range("B7").Copy()
$wsMeta.Range("B7").PasteSpecial(-4163)
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"
$wsMeta.Range("B12").Value = "ValueSet that indicates which diagnostic evidence base was used for its realization"

# Include sheet updates
$wsInclude.Range("B6").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/DiagnosisCS"
